# Atualização de bases das ligas, do dia: 28-05-2024 às 07:50
#
# The underlying data correction swaps the two match records that were
# mixed up for rows 67/68 (id 65/66, match ids 7423701/7423702) and for
# rows 125/126 (id 123/124, match ids 8039382/8039381) on the
# "Germany Verbandsliga" sheet: the HomeTeam/AwayTeam, score and odds
# columns belonging to each fixture were swapped between the two rows.
#
# We apply this by writing the corrected values directly into each
# affected cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 67 (id 65, match 7423701 -> 7423702) ---
$ws.Range("B67").Value = 7423702
$ws.Range("E67").Value = "SG RotWeiss Thalheim"
$ws.Range("F67").Value = "SV Fortuna Magdeburg"
$ws.Range("H67").Value = 0
$ws.Range("K67").Value = "D"
$ws.Range("L67").Value = 1.8
$ws.Range("N67").Value = 3
$ws.Range("O67").Value = 1.8
$ws.Range("Q67").Value = 3
$ws.Range("R67").Value = -0.5
$ws.Range("U67").Value = 3.5
$ws.Range("Y67").Value = 3.5
$ws.Range("Z67").Value = -1

# --- Row 68 (id 66, match 7423702 -> 7423701) ---
$ws.Range("B68").Value = 7423701
$ws.Range("E68").Value = "BSV HalleAmmendorf"
$ws.Range("F68").Value = "VfB Sangerhausen"
$ws.Range("H68").Value = 2
$ws.Range("K68").Value = "A"
$ws.Range("L68").Value = 2
$ws.Range("N68").Value = 2.55
$ws.Range("O68").Value = 2
$ws.Range("Q68").Value = 2.6
$ws.Range("R68").Value = -0.25
$ws.Range("U68").Value = 3.25
$ws.Range("Y68").Value = -1
$ws.Range("Z68").Value = 1.6

# --- Row 125 (id 123, match 8039382 -> 8039381) ---
$ws.Range("B125").Value = 8039381
$ws.Range("E125").Value = "SG Union Klosterfelde"
$ws.Range("F125").Value = "SV 1908 GW Ahrensfelde"
$ws.Range("G125").Value = 1
$ws.Range("H125").Value = 3
$ws.Range("I125").Value = 1
$ws.Range("J125").Value = 3
$ws.Range("L125").Value = 3.25
$ws.Range("N125").Value = 1.833
$ws.Range("O125").Value = 3.25
$ws.Range("Q125").Value = 1.833
$ws.Range("R125").Value = 0.5
$ws.Range("S125").Value = 1.925
$ws.Range("T125").Value = 1.875
$ws.Range("U125").Value = 3
$ws.Range("V125").Value = 1.825
$ws.Range("W125").Value = 1.975
$ws.Range("Z125").Value = 0.833
$ws.Range("AB125").Value = 0.875
$ws.Range("AC125").Value = 0.825
$ws.Range("AD125").Value = -1

# --- Row 126 (id 124, match 8039381 -> 8039382) ---
$ws.Range("B126").Value = 8039382
$ws.Range("E126").Value = "FC Burgsolms"
$ws.Range("F126").Value = "TSV Steinbach II"
$ws.Range("G126").Value = 0
$ws.Range("H126").Value = 4
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 2
$ws.Range("L126").Value = 15
$ws.Range("N126").Value = 1.111
$ws.Range("O126").Value = 15
$ws.Range("Q126").Value = 1.111
$ws.Range("R126").Value = 2.75
$ws.Range("S126").Value = 1.9
$ws.Range("T126").Value = 1.9
$ws.Range("U126").Value = 4
$ws.Range("V126").Value = 1.9
$ws.Range("W126").Value = 1.9
$ws.Range("Z126").Value = 0.111
$ws.Range("AB126").Value = 0.8999999999999999
$ws.Range("AC126").Value = 0
$ws.Range("AD126").Value = 0
